$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-removed trailing data rows (old rows 11-13, the Resolving-Mac -> FAPs/MuSCs/Resolving-Mac combos)
$ws.Rows("11:13").Delete()

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf13"
$ws.Cells.Item(2,3).Value = "Scn8a"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.092277
$ws.Cells.Item(2,8).Value = 0.276831
$ws.Cells.Item(2,9).Value = 0.03444274323645406
$ws.Cells.Item(2,10).Value = 0.03444274323645406
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.2858103333333333
$ws.Cells.Item(2,14).Value = 0.857431
$ws.Cells.Item(2,15).Value = 0.1287894172391936
$ws.Cells.Item(2,16).Value = 0.1287894172391936
$ws.Cells.Item(2,17).Value = 0.026373720129
$ws.Cells.Item(2,18).Value = 0.237363481161
$ws.Cells.Item(2,19).Value = 0.004435860829542095
$ws.Cells.Item(2,20).Value = 0.004435860829542095

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf13"
$ws.Cells.Item(3,3).Value = "Scn8a"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.092277
$ws.Cells.Item(3,8).Value = 0.276831
$ws.Cells.Item(3,9).Value = 0.03444274323645406
$ws.Cells.Item(3,10).Value = 0.03444274323645406
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.369620666666667
$ws.Cells.Item(3,14).Value = 4.108862
$ws.Cells.Item(3,15).Value = 0.6171667953412782
$ws.Cells.Item(3,16).Value = 0.6171667953412782
$ws.Cells.Item(3,17).Value = 0.126384486258
$ws.Cells.Item(3,18).Value = 1.137460376322
$ws.Cells.Item(3,19).Value = 0.02125691746600483
$ws.Cells.Item(3,20).Value = 0.02125691746600483

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf13"
$ws.Cells.Item(4,3).Value = "Scn8a"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.092277
$ws.Cells.Item(4,8).Value = 0.276831
$ws.Cells.Item(4,9).Value = 0.03444274323645406
$ws.Cells.Item(4,10).Value = 0.03444274323645406
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.5637756666666666
$ws.Cells.Item(4,14).Value = 1.691327
$ws.Cells.Item(4,15).Value = 0.2540437874195283
$ws.Cells.Item(4,16).Value = 0.2540437874195283
$ws.Cells.Item(4,17).Value = 0.05202352719299999
$ws.Cells.Item(4,18).Value = 0.4682117447369999
$ws.Cells.Item(4,19).Value = 0.008749964940907129
$ws.Cells.Item(4,20).Value = 0.008749964940907129

# Row 5
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Fgf13"
$ws.Cells.Item(5,3).Value = "Scn8a"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.276331666666667
$ws.Cells.Item(5,8).Value = 6.828995
$ws.Cells.Item(5,9).Value = 0.8496495022162568
$ws.Cells.Item(5,10).Value = 0.8496495022162568
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.2858103333333333
$ws.Cells.Item(5,14).Value = 0.857431
$ws.Cells.Item(5,15).Value = 0.1287894172391936
$ws.Cells.Item(5,16).Value = 0.1287894172391936
$ws.Cells.Item(5,17).Value = 0.6505991124272222
$ws.Cells.Item(5,18).Value = 5.855392011845001
$ws.Cells.Item(5,19).Value = 0.1094258642480026
$ws.Cells.Item(5,20).Value = 0.1094258642480026

# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Fgf13"
$ws.Cells.Item(6,3).Value = "Scn8a"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.276331666666667
$ws.Cells.Item(6,8).Value = 6.828995
$ws.Cells.Item(6,9).Value = 0.8496495022162568
$ws.Cells.Item(6,10).Value = 0.8496495022162568
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.369620666666667
$ws.Cells.Item(6,14).Value = 4.108862
$ws.Cells.Item(6,15).Value = 0.6171667953412782
$ws.Cells.Item(6,16).Value = 0.6171667953412782
$ws.Cells.Item(6,17).Value = 3.117710894854445
$ws.Cells.Item(6,18).Value = 28.05939805369
$ws.Cells.Item(6,19).Value = 0.5243754604461194
$ws.Cells.Item(6,20).Value = 0.5243754604461194

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Fgf13"
$ws.Cells.Item(7,3).Value = "Scn8a"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.276331666666667
$ws.Cells.Item(7,8).Value = 6.828995
$ws.Cells.Item(7,9).Value = 0.8496495022162568
$ws.Cells.Item(7,10).Value = 0.8496495022162568
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.5637756666666666
$ws.Cells.Item(7,14).Value = 1.691327
$ws.Cells.Item(7,15).Value = 0.2540437874195283
$ws.Cells.Item(7,16).Value = 0.2540437874195283
$ws.Cells.Item(7,17).Value = 1.283340402929444
$ws.Cells.Item(7,18).Value = 11.550063626365
$ws.Cells.Item(7,19).Value = 0.2158481775221348
$ws.Cells.Item(7,20).Value = 0.2158481775221348

# Row 8
$ws.Cells.Item(8,1).Value = "Resolving-Mac"
$ws.Cells.Item(8,2).Value = "Fgf13"
$ws.Cells.Item(8,3).Value = "Scn8a"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.3105333333333333
$ws.Cells.Item(8,8).Value = 0.9316
$ws.Cells.Item(8,9).Value = 0.1159077545472891
$ws.Cells.Item(8,10).Value = 0.1159077545472891
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.2858103333333333
$ws.Cells.Item(8,14).Value = 0.857431
$ws.Cells.Item(8,15).Value = 0.1287894172391936
$ws.Cells.Item(8,16).Value = 0.1287894172391936
$ws.Cells.Item(8,17).Value = 0.08875363551111111
$ws.Cells.Item(8,18).Value = 0.7987827196
$ws.Cells.Item(8,19).Value = 0.01492769216164886
$ws.Cells.Item(8,20).Value = 0.01492769216164886

# Row 9
$ws.Cells.Item(9,1).Value = "Resolving-Mac"
$ws.Cells.Item(9,2).Value = "Fgf13"
$ws.Cells.Item(9,3).Value = "Scn8a"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.3105333333333333
$ws.Cells.Item(9,8).Value = 0.9316
$ws.Cells.Item(9,9).Value = 0.1159077545472891
$ws.Cells.Item(9,10).Value = 0.1159077545472891
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.369620666666667
$ws.Cells.Item(9,14).Value = 4.108862
$ws.Cells.Item(9,15).Value = 0.6171667953412782
$ws.Cells.Item(9,16).Value = 0.6171667953412782
$ws.Cells.Item(9,17).Value = 0.4253128710222223
$ws.Cells.Item(9,18).Value = 3.8278158392
$ws.Cells.Item(9,19).Value = 0.0715344174291539
$ws.Cells.Item(9,20).Value = 0.07153441742915391

# Row 10
$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Fgf13"
$ws.Cells.Item(10,3).Value = "Scn8a"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.3105333333333333
$ws.Cells.Item(10,8).Value = 0.9316
$ws.Cells.Item(10,9).Value = 0.1159077545472891
$ws.Cells.Item(10,10).Value = 0.1159077545472891
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.5637756666666666
$ws.Cells.Item(10,14).Value = 1.691327
$ws.Cells.Item(10,15).Value = 0.2540437874195283
$ws.Cells.Item(10,16).Value = 0.2540437874195283
$ws.Cells.Item(10,17).Value = 0.1750711370222222
$ws.Cells.Item(10,18).Value = 1.5756402332
$ws.Cells.Item(10,19).Value = 0.02944564495648638
$ws.Cells.Item(10,20).Value = 0.02944564495648638
